# edit.ps1
# Applies the betting-odds value updates described in the commit diff
# for "Jogos_da_Semana_FlashScore_FULL_2024-11-22.xlsx" (rows 3-7, 16, 17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("AF3").Value = 41
$ws.Range("AJ3").Value = 19
$ws.Range("AS3").Value = 101
$ws.Range("AT3").Value = 3.4
$ws.Range("AU3").Value = 8
$ws.Range("BB3").Value = 101
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 6.25
$ws.Range("L3").Value = 5.5
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.4
$ws.Range("U3").Value = 1.73
$ws.Range("V3").Value = 2

# Row 4
$ws.Range("AJ4").Value = 13
$ws.Range("AO4").Value = 13
$ws.Range("AS4").Value = 301
$ws.Range("G4").Value = 2.1
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 3.75
$ws.Range("J4").Value = 3
$ws.Range("L4").Value = 4.5

# Row 5
$ws.Range("AE5").Value = 19
$ws.Range("AJ5").Value = 19
$ws.Range("G5").Value = 1.62
$ws.Range("I5").Value = 6
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.7
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67

# Row 6
$ws.Range("AD6").Value = 7.5
$ws.Range("AK6").Value = 81
$ws.Range("BB6").Value = 201
$ws.Range("G6").Value = 1.53
$ws.Range("H6").Value = 3.75
$ws.Range("J6").Value = 2.1
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 10

# Row 7
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 9.5
$ws.Range("AH7").Value = 9.5
$ws.Range("AI7").Value = 17
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 34
$ws.Range("AM7").Value = 34
$ws.Range("AO7").Value = 12
$ws.Range("AP7").Value = 23
$ws.Range("AQ7").Value = 41
$ws.Range("AT7").Value = 2.63
$ws.Range("AY7").Value = 19
$ws.Range("BB7").Value = 81
$ws.Range("G7").Value = 2.1
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 3.3
$ws.Range("J7").Value = 2.88
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 2.05
$ws.Range("R7").Value = 1.8
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.63
$ws.Range("X7").Value = 10
$ws.Range("Z7").Value = 19

# Row 16
$ws.Range("M16").Value = 1.01
$ws.Range("N16").Value = 16.5
$ws.Range("O16").Value = $null
$ws.Range("P16").Value = $null
$ws.Range("Q16").Value = 1.19
$ws.Range("R16").Value = 3.6
$ws.Range("S16").Value = 1.13
$ws.Range("T16").Value = 5.1

# Row 17
$ws.Range("AB17").Value = 26
$ws.Range("AE17").Value = 22
$ws.Range("AF17").Value = 90
$ws.Range("AG17").Value = 700
$ws.Range("AH17").Value = 30
$ws.Range("AI17").Value = 80
$ws.Range("AJ17").Value = 30
$ws.Range("AK17").Value = 300
$ws.Range("AL17").Value = 120
$ws.Range("AP17").Value = 14
$ws.Range("AU17").Value = 8.5
$ws.Range("AX17").Value = 10.5
$ws.Range("AY17").Value = 55
$ws.Range("BB17").Value = 300
$ws.Range("H17").Value = 5.3
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 1.65
$ws.Range("K17").Value = 2.67
$ws.Range("L17").Value = 7.8
$ws.Range("O17").Value = 1.16
$ws.Range("P17").Value = 4.55
$ws.Range("Q17").Value = 1.5
$ws.Range("R17").Value = 2.4
$ws.Range("U17").Value = 1.93
$ws.Range("V17").Value = 1.78
$ws.Range("Y17").Value = 8.75
